$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 20 (pushes existing rows 20+ down by two)
$ws.Rows("20:21").Insert()

# The inserted rows don't automatically carry the same cell borders/format
# as the surrounding table, so copy the formatting down from row 19.
$ws.Range("A19:F19").Copy()
$ws.Range("A20:F21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the two new rows with the completed TP5 tasks
$ws.Range("A20").Value = "Améliore le mouvement des caisses"
$ws.Range("B20").Value = "30min"
$ws.Range("C20").Value = "TP5"
$ws.Range("F20").Value = "OK"

$ws.Range("A21").Value = "Termine une map compléte"
$ws.Range("B21").Value = "2h"
$ws.Range("C21").Value = "TP5-TP6"
$ws.Range("F21").Value = "OK"

# Update the active selection to match the author's final cursor position
$ws.Range("A25").Select()
